$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.165568306009249
$ws.Cells.Item(2, 3).Value = 0.0334234178857109
$ws.Cells.Item(2, 4).Value = 0.3090660498759803
$ws.Cells.Item(2, 5).Value = 0.08600710581194093
$ws.Cells.Item(2, 6).Value = 4.618477121752477
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 11).Value = 0.6763136704455803
$ws.Cells.Item(2, 12).Value = 0.2225196855931131

$ws.Cells.Item(3, 2).Value = 1.145630774578621
$ws.Cells.Item(3, 3).Value = 0.03160631916467338
$ws.Cells.Item(3, 4).Value = 0.2967283011829522
$ws.Cells.Item(3, 5).Value = 0.08467592048555161
$ws.Cells.Item(3, 6).Value = 4.405935134467256
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 11).Value = 0.6589400703126387
$ws.Cells.Item(3, 12).Value = 0.2168610581895081

$ws.Cells.Item(4, 2).Value = 1.134448267733461
$ws.Cells.Item(4, 3).Value = 0.03046620734637884
$ws.Cells.Item(4, 4).Value = 0.2890757357964588
$ws.Cells.Item(4, 5).Value = 0.08388671296320993
$ws.Cells.Item(4, 6).Value = 4.27568200701856
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 11).Value = 0.6490472067142861
$ws.Cells.Item(4, 12).Value = 0.2135342625218755

$ws.Cells.Item(5, 2).Value = 1.130157057674921
$ws.Cells.Item(5, 3).Value = 0.02999531420208257
$ws.Cells.Item(5, 4).Value = 0.2859374313792102
$ws.Cells.Item(5, 5).Value = 0.08357216748360585
$ws.Cells.Item(5, 6).Value = 4.222660805723308
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 11).Value = 0.6452098223988401
$ws.Cells.Item(5, 12).Value = 0.2122155782463295

$ws.Cells.Item(6, 2).Value = 1.129460539230024
$ws.Cells.Item(6, 3).Value = 0.0299167385449195
$ws.Cells.Item(6, 4).Value = 0.2854151087299499
$ws.Cells.Item(6, 5).Value = 0.08352036362281723
$ws.Cells.Item(6, 6).Value = 4.213860048901381
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 11).Value = 0.6445843257182844
$ws.Cells.Item(6, 12).Value = 0.2119988440229292

$ws.Cells.Item(7, 2).Value = 1.134389319740308
$ws.Cells.Item(7, 3).Value = 0.03045988237635555
$ws.Cells.Item(7, 4).Value = 0.2890334923381914
$ws.Cells.Item(7, 5).Value = 0.08388244231311859
$ws.Cells.Item(7, 6).Value = 4.274966715616443
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 11).Value = 0.6489946697320903
$ws.Cells.Item(7, 12).Value = 0.213516328565774

$ws.Cells.Item(8, 2).Value = 1.158473629041339
$ws.Cells.Item(8, 3).Value = 0.03280185618442033
$ws.Cells.Item(8, 4).Value = 0.304827745216528
$ws.Cells.Item(8, 5).Value = 0.0855422621095272
$ws.Cells.Item(8, 6).Value = 4.545138293844047
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 11).Value = 0.6701620266318002
$ws.Cells.Item(8, 12).Value = 0.2205378824116906

$ws.Cells.Item(9, 2).Value = 1.214139903148293
$ws.Cells.Item(9, 3).Value = 0.03720773664952404
$ws.Cells.Item(9, 4).Value = 0.3352084593906426
$ws.Cells.Item(9, 5).Value = 0.08902139005117249
$ws.Cells.Item(9, 6).Value = 5.077144509861768
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 11).Value = 0.7178550510727177
$ws.Cells.Item(9, 12).Value = 0.2354845402657446

$ws.Cells.Item(10, 2).Value = 1.260234718278213
$ws.Cells.Item(10, 3).Value = 0.04034064433363227
$ws.Cells.Item(10, 4).Value = 0.3571981080038995
$ws.Cells.Item(10, 5).Value = 0.09171581356577718
$ws.Cells.Item(10, 6).Value = 5.469705866913387
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 11).Value = 0.7567237436143728
$ws.Cells.Item(10, 12).Value = 0.2471937039238838

$ws.Cells.Item(11, 2).Value = 1.282345237153834
$ws.Cells.Item(11, 3).Value = 0.04174540326025777
$ws.Cells.Item(11, 4).Value = 0.3671363910970911
$ws.Cells.Item(11, 5).Value = 0.09297198724314626
$ws.Cells.Item(11, 6).Value = 5.648740027691247
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 11).Value = 0.7752508165879703
$ws.Cells.Item(11, 12).Value = 0.2526808745670053

$ws.Cells.Item(12, 2).Value = 1.290882992621647
$ws.Cells.Item(12, 3).Value = 0.0422745836962477
$ws.Cells.Item(12, 4).Value = 0.3708909178179738
$ws.Cells.Item(12, 5).Value = 0.09345207302573044
$ws.Cells.Item(12, 6).Value = 5.716606895914367
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 11).Value = 0.7823891224854265
$ws.Cells.Item(12, 12).Value = 0.2547819814463281

$ws.Cells.Item(13, 2).Value = 1.289036883633713
$ws.Cells.Item(13, 3).Value = 0.04216073583808111
$ws.Cells.Item(13, 4).Value = 0.3700827013004471
$ws.Cells.Item(13, 5).Value = 0.09334848195986467
$ws.Cells.Item(13, 6).Value = 5.701987335881995
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 11).Value = 0.7808462972132588
$ws.Cells.Item(13, 12).Value = 0.2543284348854655

$ws.Cells.Item(14, 2).Value = 1.283044332293031
$ws.Cells.Item(14, 3).Value = 0.04178899384593393
$ws.Cells.Item(14, 4).Value = 0.3674454537592169
$ws.Cells.Item(14, 5).Value = 0.09301139589449292
$ws.Cells.Item(14, 6).Value = 5.654322033977451
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 11).Value = 0.775835628760575
$ws.Cells.Item(14, 12).Value = 0.2528532674354125

$ws.Cells.Item(15, 2).Value = 1.279395231182775
$ws.Cells.Item(15, 3).Value = 0.04156093515011605
$ws.Cells.Item(15, 4).Value = 0.3658289197512943
$ws.Cells.Item(15, 5).Value = 0.09280549453920628
$ws.Cells.Item(15, 6).Value = 5.625134991443588
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 11).Value = 0.7727824312115956
$ws.Cells.Item(15, 12).Value = 0.251952714760165

$ws.Cells.Item(16, 2).Value = 1.258812786651845
$ws.Cells.Item(16, 3).Value = 0.04024844323686239
$ws.Cells.Item(16, 4).Value = 0.3565473498899223
$ws.Cells.Item(16, 5).Value = 0.09163433434172674
$ws.Cells.Item(16, 6).Value = 5.458015160744338
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 11).Value = 0.7555300479886
$ws.Cells.Item(16, 12).Value = 0.2468383481203347

$ws.Cells.Item(17, 2).Value = 1.246479115466002
$ws.Cells.Item(17, 3).Value = 0.03943815943855355
$ws.Cells.Item(17, 4).Value = 0.350837099378623
$ws.Cells.Item(17, 5).Value = 0.09092368210847823
$ws.Cells.Item(17, 6).Value = 5.355612851994124
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 11).Value = 0.7451634744682565
$ws.Cells.Item(17, 12).Value = 0.2437420841302469

$ws.Cells.Item(18, 2).Value = 1.239492532319787
$ws.Cells.Item(18, 3).Value = 0.03897017399986424
$ws.Cells.Item(18, 4).Value = 0.3475465768217134
$ws.Cells.Item(18, 5).Value = 0.09051780025870571
$ws.Cells.Item(18, 6).Value = 5.296756376690837
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 11).Value = 0.7392804225294185
$ws.Cells.Item(18, 12).Value = 0.2419763041688867

$ws.Cells.Item(19, 2).Value = 1.237145420554185
$ws.Cells.Item(19, 3).Value = 0.03881138542279672
$ws.Cells.Item(19, 4).Value = 0.3464313946491444
$ws.Cells.Item(19, 5).Value = 0.09038086739766271
$ws.Cells.Item(19, 6).Value = 5.276835756767696
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 11).Value = 0.7373021532995381
$ws.Cells.Item(19, 12).Value = 0.2413810322732672

$ws.Cells.Item(20, 2).Value = 1.247780932725789
$ws.Cells.Item(20, 3).Value = 0.03952461426116116
$ws.Cells.Item(20, 4).Value = 0.3514455972519102
$ws.Cells.Item(20, 5).Value = 0.09099903544255028
$ws.Cells.Item(20, 6).Value = 5.366509305552256
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 11).Value = 0.7462587765673732
$ws.Cells.Item(20, 12).Value = 0.2440701219824177

$ws.Cells.Item(21, 2).Value = 1.28480000583906
$ws.Cells.Item(21, 3).Value = 0.0418982573732265
$ws.Cells.Item(21, 4).Value = 0.3682203148375436
$ws.Cells.Item(21, 5).Value = 0.09311028667224619
$ws.Cells.Item(21, 6).Value = 5.668320535046121
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 11).Value = 0.7773040513272917
$ws.Cells.Item(21, 12).Value = 0.2532859281766662

$ws.Cells.Item(22, 2).Value = 1.309956138354551
$ws.Cells.Item(22, 3).Value = 0.04343349468096847
$ws.Cells.Item(22, 4).Value = 0.3791319824927655
$ws.Cells.Item(22, 5).Value = 0.09451576762985425
$ws.Cells.Item(22, 6).Value = 5.865984989641731
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 11).Value = 0.7983084318170199
$ws.Cells.Item(22, 12).Value = 0.2594444852623212

$ws.Cells.Item(23, 2).Value = 1.296441541160618
$ws.Cells.Item(23, 3).Value = 0.04261552578533667
$ws.Cells.Item(23, 4).Value = 0.3733127902918625
$ws.Cells.Item(23, 5).Value = 0.09376328201953044
$ws.Cells.Item(23, 6).Value = 5.760448276784757
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 11).Value = 0.7870323116867439
$ws.Cells.Item(23, 12).Value = 0.2561451029434778

$ws.Cells.Item(24, 2).Value = 1.247192056615802
$ws.Cells.Item(24, 3).Value = 0.03948553472231708
$ws.Cells.Item(24, 4).Value = 0.3511705191425278
$ws.Cells.Item(24, 5).Value = 0.09096495984908515
$ws.Cells.Item(24, 6).Value = 5.361582969508987
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 11).Value = 0.7457633512516395
$ws.Cells.Item(24, 12).Value = 0.2439217715281927

$ws.Cells.Item(25, 2).Value = 1.198172025980995
$ws.Cells.Item(25, 3).Value = 0.03603481746144865
$ws.Cells.Item(25, 4).Value = 0.3270499636984283
$ws.Cells.Item(25, 5).Value = 0.08805602403618451
$ws.Cells.Item(25, 6).Value = 4.932950692271561
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 11).Value = 0.704284794220456
$ws.Cells.Item(25, 12).Value = 0.2313139945974712

